# Generate Report for Handoff
# Update Priority and Latest Handoff Datetime for the "Ready for handoff"
# rows (4-7) on both the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-09-01 06:34:44"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-09-01 06:34:49"
